$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the "Meta description" paragraph that currently sits right
#    after the H1 title paragraph.
# ------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Meta description:*") {
        $metaPara = $candidate
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) At the tail of the document, turn the single italic "feature
#    image prompt" paragraph into two paragraphs:
#      - a new bold paragraph carrying the page title text, and
#      - the same italic paragraph, retexted to the meta-description
#        copy that used to live at the top of the document.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)

# The document-editing host has a quirk where InsertXML-ing a range
# that reaches all the way to the end of the story silently drops the
# trailing content and leaves a stray empty paragraph. Sidestep that
# by first appending a scratch paragraph after the last paragraph, so
# the range we operate on no longer touches the very end of the
# story; we delete the scratch paragraph again afterwards.
$lastPara.Range.InsertParagraphAfter()

$target = $d.Paragraphs.Item($n)
$targetRange = $target.Range

$replacementXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Charlie Chance In Hell To Pay for Free | Slot Game Review</w:t></w:r></w:p>' +
  '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Charlie Chance In Hell To Pay, a unique and moderately volatile slot game with five bonus features. Play for free and try your luck!</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($replacementXml)

# Remove the scratch paragraph left dangling at the very end of the story.
$scratch = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratch.Range.Delete()
